# Generate Report for Handback
# Adds a new handback row (GUID b6fac65e-5138-4ca2-8939-d7d4628f4b86) to all
# three worksheets: Overview (row 4, A:C) and zh-cn / de-de (row 4, A:L).

$wb = $excel.ActiveWorkbook

$guid   = "b6fac65e-5138-4ca2-8939-d7d4628f4b86"
$mdName = "$guid.md"
$xlfBaseZh = "$guid.802869d445d04ef651ae7cd343bc7a3450df34e5.zh-cn.xlf"
$xlfBaseDe = "$guid.802869d445d04ef651ae7cd343bc7a3450df34e5.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"

# Adds a hyperlink to a cell and paints it with the same look as the
# workbook's existing custom "HyperLink" cell style (single underline,
# RGB 6495ED -- 15570276 once packed into an OLE/BGR color value).
function Set-HL($ws, $cellRef, $address, $display) {
    $rng = $ws.Range($cellRef)
    $rng.Value = $display
    $ws.Hyperlinks.Add($rng, $address, "", "", $display) | Out-Null
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Underline = 2
    $rng.Font.Color = 15570276
}

# Stamps a literal text timestamp (not a real Excel date serial -- the
# source workbook stores these as text) with the existing datetime
# number format so the cell style matches rows 2/3.
function Set-DateText($ws, $cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = $text
    $rng.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ----------------------------------------------------------------------
# Sheet 1: Overview  (A4:C4)
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-HL $wsOverview "A4" "https://github.com/OpenLocalizationTest/oltest/blob/3db0623068ea2f59aaa8a65b2d103e3dd244dd41/e2e/$mdName" $mdName

$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ----------------------------------------------------------------------
# Sheet 2: zh-cn  (A4:J4, columns K/L left blank like rows 2-3)
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HL $wsZh "A4" "https://github.com/OpenLocalizationTest/oltest/blob/3db0623068ea2f59aaa8a65b2d103e3dd244dd41/e2e/$mdName" $mdName

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $statusInSync

Set-HL $wsZh "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f073ff4d7a9aee02b6e7dd1c0bda6af8ceab021/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$xlfBaseZh" $xlfBaseZh

Set-DateText $wsZh "E4" "2016-03-19 03:28:52"

Set-HL $wsZh "F4" "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/01cfd7a47380f373a9240216b8261b6960008b76/e2e/$mdName" $mdName

Set-HL $wsZh "G4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eca8ca2cc381924d00eab307653762e0688d0950/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$xlfBaseZh" $xlfBaseZh

Set-DateText $wsZh "H4" "2016-03-19 03:29:31"

$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = "Include"

# ----------------------------------------------------------------------
# Sheet 3: de-de  (A4:J4, columns K/L left blank like rows 2-3)
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-HL $wsDe "A4" "https://github.com/OpenLocalizationTest/oltest/blob/3db0623068ea2f59aaa8a65b2d103e3dd244dd41/e2e/$mdName" $mdName

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $statusInSync

Set-HL $wsDe "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d4f43fcccad04e02d4ba58c25eeb87f69f884d31/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$xlfBaseDe" $xlfBaseDe

Set-DateText $wsDe "E4" "2016-03-19 03:29:02"

Set-HL $wsDe "F4" "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/2885ead5c7db613ff4558c46d7c399f40a9769d2/e2e/$mdName" $mdName

Set-HL $wsDe "G4" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8fea263f8e83c5e622bdb6fb0e4fa2f7d97f6ac6/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$xlfBaseDe" $xlfBaseDe

Set-DateText $wsDe "H4" "2016-03-19 03:29:44"

$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = "Include"
